# Update "想去人数" (F column) values across the four worksheets to match
# the latest scrape output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 2678
$ws.Range("F10").Value = 5941
$ws.Range("F16").Value = 4888
$ws.Range("F20").Value = 2497
$ws.Range("F22").Value = 485
$ws.Range("F25").Value = 268
$ws.Range("F30").Value = 383
$ws.Range("F37").Value = 1447
$ws.Range("F41").Value = 223
$ws.Range("F42").Value = 1705
$ws.Range("F43").Value = 2485
$ws.Range("F45").Value = 100
$ws.Range("F49").Value = 82

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value  = 399
$ws.Range("F22").Value = 317

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value  = 1410
$ws.Range("F9").Value  = 1794
$ws.Range("F10").Value = 2362
$ws.Range("F11").Value = 774
$ws.Range("F12").Value = 657

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 2678
$ws.Range("F8").Value  = 1410
$ws.Range("F10").Value = 2362
$ws.Range("F11").Value = 5941
$ws.Range("F12").Value = 774
$ws.Range("F16").Value = 4888
$ws.Range("F17").Value = 2497
$ws.Range("F19").Value = 485
$ws.Range("F26").Value = 383
$ws.Range("F33").Value = 1447
$ws.Range("F39").Value = 317
$ws.Range("F40").Value = 223
$ws.Range("F42").Value = 1705
$ws.Range("F43").Value = 2485
$ws.Range("F44").Value = 100

$wb.Save()
